$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Checks" sheet: null_check / count_check / column_check now default to "Yes"
# ---------------------------------------------------------------------------
$wsChecks = $wb.Worksheets.Item("Checks")
$wsChecks.Range("B2").Value = "Yes"
$wsChecks.Range("B3").Value = "Yes"
$wsChecks.Range("B4").Value = "Yes"

# ---------------------------------------------------------------------------
# "UploadDataToDB" sheet: add etl/table_name/table_creation_sql_query rows
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("UploadDataToDB")

# Make room for the two new rows (table_name / table_creation_sql_query)
# right after the existing "database" row - this pushes the csv/json/excel
# drop-down list (rows 28-30) down to rows 29-31.
$ws.Rows.Item(9).Insert()

# Populate the new cells. The order in which *new* distinct string values
# are first written controls the shared-string table ordering, so row 10 is
# populated before row 9.
$ws.Range("A10").Value = "table_creation_sql_query"
$ws.Range("A9").Value = "table_name"
$ws.Range("B8").Value = "etl"
$ws.Range("B9").Value = "insertedData"
$ws.Range("B10").Value = "CREATE TABLE insertedData (Id INT , empid INT, name VARCHAR(255), designation VARCHAR(255))"

# The SQL cell wraps and the row grows to fit it.
$ws.Range("A10").WrapText = $true
$ws.Rows.Item(10).RowHeight = 43.5

# Widen the two columns so the new labels / SQL text are readable.
$ws.Columns.Item(1).ColumnWidth = 15
$ws.Columns.Item(2).ColumnWidth = 24.5

# The drop-down validation list moved down by one row.
$ws.Range("B1").Validation.Delete()
$ws.Range("B1").Validation.Add(3, 1, 1, "=`$A`$29:`$A`$32")

# ---------------------------------------------------------------------------
# Selection bookkeeping - UploadDataToDB's own cursor moves to B10, but the
# workbook's active tab stays on "Checks" (also now parked at B10).
# ---------------------------------------------------------------------------
$ws.Activate()
[void]$ws.Range("B10").Select()

$wsChecks.Activate()
[void]$wsChecks.Range("B10").Select()
